# Refactorizacion general - actualizacion de ruta (v2.0-alpha)
# - Se elimina la fila 3 (ACHS Osorno), lo que desplaza hacia arriba el resto
#   de las filas de datos (4..18 -> 3..17) y actualiza el rango de uso.
# - Se agregan dos clientes nuevos al final de la ruta (filas 18 y 19).
# - Se corrige el formato de la fila 64 (antes mezclaba estilos).
# - La fila 138 queda vacia (arrastrada por el corrimiento de filas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Eliminar la fila 3 completa; Excel recorre hacia arriba el resto de filas.
$ws.Rows(3).Delete()

# 2. Completar la nueva fila 18 con los datos del cliente "Maria Jose Rodriguez".
#    Las columnas A y B deben quedar como valores numericos (aunque la celda
#    tenga aplicado un formato de texto), por lo que se cambia momentaneamente
#    el formato a General para el ingreso y luego se repone el formato de texto.
$ws.Range("A18").NumberFormat = "General"
$ws.Range("A18").Value2 = 20250318
$ws.Range("A18").NumberFormat = "@"

$ws.Range("B18").NumberFormat = "General"
$ws.Range("B18").Value2 = 16
$ws.Range("B18").NumberFormat = "@"

$ws.Range("C18").Value2 = "30"
$ws.Range("D18").Value2 = "17.673.326-8"
$ws.Range("E18").Value2 = "Maria José Rodriguez"
$ws.Range("F18").Value2 = "Colaco s/n km 3, parcela 9"
$ws.Range("G18").Value2 = "Calbuco"
$ws.Range("H18").Value2 = "972861950"
$ws.Range("I18").Value2 = "Cliente test"
$ws.Range("J18").Value2 = "1002"

# 3. Completar la nueva fila 19 con los datos del cliente "Isaias Beroiza Mora".
$ws.Range("A19").NumberFormat = "General"
$ws.Range("A19").Value2 = 20250318
$ws.Range("A19").NumberFormat = "@"

$ws.Range("B19").NumberFormat = "General"
$ws.Range("B19").Value2 = 17
$ws.Range("B19").NumberFormat = "@"

$ws.Range("C19").Value2 = "30"
$ws.Range("D19").Value2 = "16.742.249-7"
$ws.Range("E19").Value2 = "Isaias Beroiza Mora"
$ws.Range("F19").Value2 = "Colaco s/n km 3, parcela 9"
$ws.Range("G19").Value2 = "Calbuco"
$ws.Range("H19").Value2 = "988809704"
$ws.Range("I19").Value2 = "Cliente test"
$ws.Range("J19").Value2 = "1001"

# 4. Corregir el formato de la fila 64 para que coincida con el resto de filas
#    vacias (columna A, G, H, I, J deben compartir el mismo estilo que B:F).
$ws.Range("A65:J65").Copy()
$ws.Range("A64:J64").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
